$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1375.7174
$ws.Range("I15").Value = 1375.7174
$ws.Range("K15").Value = 4127.1522
$ws.Range("M15").Value = -3958.1522

$ws.Range("H62").Value = 4555.5557
$ws.Range("I62").Value = 3833.3333
$ws.Range("K62").Value = 3833.3333
$ws.Range("M62").Value = -3209.3333

$ws.Range("H65").Value = 4555.5557
$ws.Range("I65").Value = 3833.3333
$ws.Range("K65").Value = 19166.6665
$ws.Range("M65").Value = -16046.6665

$ws.Range("H69").Value = 55568450
$ws.Range("I69").Value = 111116900
$ws.Range("J69").Value = 20000
$ws.Range("K69").Value = 333350700
$ws.Range("L69").Value = 60000
$ws.Range("M69").Value = -333349826
$ws.Range("N69").Value = -61748

$ws.Range("H72").Value = 55568450
$ws.Range("I72").Value = 111116900
$ws.Range("J72").Value = 20000
$ws.Range("K72").Value = 1000052100
$ws.Range("L72").Value = 180000
$ws.Range("M72").Value = -1000047732
$ws.Range("N72").Value = -188736

$ws.Range("H86").Value = 111113096
$ws.Range("I86").Value = 166669150
$ws.Range("K86").Value = 166669150
$ws.Range("M86").Value = -166668027

$ws.Range("H89").Value = 111113096
$ws.Range("I89").Value = 166669150
$ws.Range("K89").Value = 833345750
$ws.Range("M89").Value = -833340134

$ws.Range("H96").Value = 5727.6924
$ws.Range("J96").Value = 13897.6
$ws.Range("L96").Value = 41692.8
$ws.Range("N96").Value = -44438.8

$ws.Range("H111").Value = 4840.905
$ws.Range("I111").Value = 3148.5
$ws.Range("K111").Value = 9445.5
$ws.Range("M111").Value = -6378.5

$ws.Range("H116").Value = 24641730
$ws.Range("I116").Value = 25761452
$ws.Range("J116").Value = 7856
$ws.Range("K116").Value = 25761452
$ws.Range("L116").Value = 7856
$ws.Range("M116").Value = -25758010
$ws.Range("N116").Value = -14740

$ws.Range("H141").Value = 1589.7778
$ws.Range("I141").Value = 1589.7778
$ws.Range("K141").Value = 4769.3334
$ws.Range("M141").Value = 410.6665999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2501.8914
$ws.Range("I134").Value = 2112.8518
$ws.Range("K134").Value = 6338.555399999999
$ws.Range("M134").Value = -3803.555399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2530.4285
$ws.Range("I16").Value = 1283.3334
$ws.Range("K16").Value = 1283.3334
$ws.Range("M16").Value = -996.3334

$ws.Range("H31").Value = 15628234
$ws.Range("I31").Value = 17546780
$ws.Range("J31").Value = 5782.143
$ws.Range("K31").Value = 17546780
$ws.Range("L31").Value = 5782.143
$ws.Range("M31").Value = -17546485
$ws.Range("N31").Value = -6372.143

$ws.Range("H34").Value = 15628234
$ws.Range("I34").Value = 17546780
$ws.Range("J34").Value = 5782.143
$ws.Range("K34").Value = 17546780
$ws.Range("L34").Value = 5782.143
$ws.Range("M34").Value = -17546578
$ws.Range("N34").Value = -6186.143

$ws.Range("H107").Value = 652.03845
$ws.Range("I107").Value = 406.18182
$ws.Range("K107").Value = 406.18182
$ws.Range("M107").Value = 1513.81818

$ws.Range("H113").Value = 2530.4285
$ws.Range("I113").Value = 1283.3334
$ws.Range("K113").Value = 1283.3334
$ws.Range("M113").Value = 886.6666

$ws.Range("H132").Value = 40405316
$ws.Range("I132").Value = 49383884
$ws.Range("J132").Value = 1765.1666
$ws.Range("K132").Value = 148151652
$ws.Range("L132").Value = 5295.4998
$ws.Range("M132").Value = -148149122
$ws.Range("N132").Value = -10355.4998

$ws.Range("H134").Value = 2754.5862
$ws.Range("I134").Value = 2285.2
$ws.Range("K134").Value = 6855.599999999999
$ws.Range("M134").Value = -4320.599999999999

$ws.Range("H135").Value = 53666.25
$ws.Range("I135").Value = 40000
$ws.Range("J135").Value = 76443.336
$ws.Range("K135").Value = 40000
$ws.Range("L135").Value = 76443.336
$ws.Range("M135").Value = -34930
$ws.Range("N135").Value = -86583.336

$ws.Range("H141").Value = 125741.76
$ws.Range("I141").Value = 86333.164
$ws.Range("J141").Value = 134499.22
$ws.Range("K141").Value = 86333.164
$ws.Range("L141").Value = 134499.22
$ws.Range("M141").Value = -81153.164
$ws.Range("N141").Value = -144859.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 4354.4
$ws.Range("I49").Value = 595.3333
$ws.Range("K49").Value = 1785.9999
$ws.Range("M49").Value = -1629.9999

$ws.Range("H95").Value = 18949.5
$ws.Range("J95").Value = 18949.5
$ws.Range("L95").Value = 56848.5
$ws.Range("N95").Value = -60966.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 47350.8
$ws.Range("J62").Value = 45584.668
$ws.Range("L62").Value = 45584.668
$ws.Range("N62").Value = -46956.668

$ws.Range("H65").Value = 47350.8
$ws.Range("J65").Value = 45584.668
$ws.Range("L65").Value = 136754.004
$ws.Range("N65").Value = -143618.004

$ws.Range("H126").Value = 8876.5
$ws.Range("I126").Value = 11790.5
$ws.Range("J126").Value = 5962.5
$ws.Range("K126").Value = 35371.5
$ws.Range("L126").Value = 17887.5
$ws.Range("M126").Value = -32901.5
$ws.Range("N126").Value = -22827.5

$ws.Range("H132").Value = 68506.664
$ws.Range("I132").Value = 81642.16
$ws.Range("J132").Value = 2829.2
$ws.Range("K132").Value = 244926.48
$ws.Range("L132").Value = 8487.599999999999
$ws.Range("M132").Value = -242396.48
$ws.Range("N132").Value = -13547.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6198.905
$ws.Range("I7").Value = 5555.9
$ws.Range("J7").Value = 6783.4546
$ws.Range("K7").Value = 5555.9
$ws.Range("L7").Value = 6783.4546
$ws.Range("M7").Value = -5443.9
$ws.Range("N7").Value = -7007.4546

$ws.Range("H22").Value = 948.6875
$ws.Range("I22").Value = 815
$ws.Range("K22").Value = 815
$ws.Range("M22").Value = -520

$ws.Range("H27").Value = 948.6875
$ws.Range("I27").Value = 815
$ws.Range("K27").Value = 815
$ws.Range("M27").Value = -708

$ws.Range("H40").Value = 20707640
$ws.Range("I40").Value = 11907766
$ws.Range("K40").Value = 11907766
$ws.Range("M40").Value = -11907630

$ws.Range("H126").Value = 6198.905
$ws.Range("I126").Value = 5555.9
$ws.Range("J126").Value = 6783.4546
$ws.Range("K126").Value = 16667.7
$ws.Range("L126").Value = 20350.3638
$ws.Range("M126").Value = -14197.7
$ws.Range("N126").Value = -25290.3638

$ws.Range("H132").Value = 2669.04
$ws.Range("I132").Value = 2638.959
$ws.Range("J132").Value = 2750.3704
$ws.Range("K132").Value = 7916.876999999999
$ws.Range("L132").Value = 8251.111199999999
$ws.Range("M132").Value = -5386.876999999999
$ws.Range("N132").Value = -13311.1112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 39333.332
$ws.Range("J110").Value = 39333.332
$ws.Range("L110").Value = 39333.332
$ws.Range("N110").Value = -47513.332

$ws.Range("H126").Value = 50002580
$ws.Range("I126").Value = 76925580
$ws.Range("J126").Value = 2712.8572
$ws.Range("K126").Value = 230776740
$ws.Range("L126").Value = 8138.571599999999
$ws.Range("M126").Value = -230774270
$ws.Range("N126").Value = -13078.5716

$ws.Range("H132").Value = 4656.382
$ws.Range("I132").Value = 1014.1739
$ws.Range("J132").Value = 7274.2188
$ws.Range("K132").Value = 3042.5217
$ws.Range("L132").Value = 21822.6564
$ws.Range("M132").Value = -512.5217000000002
$ws.Range("N132").Value = -26882.6564

$ws.Range("H136").Value = 4359.6055
$ws.Range("I136").Value = 2784.4443
$ws.Range("K136").Value = 8353.332900000001
$ws.Range("M136").Value = -5803.332900000001
